$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy row 113's formatting for columns A (bold/border id style) and E
# (custom date number format) down into rows 114 and 115, matching the
# formatting already used throughout the rest of the table.
$ws.Range("A113").Copy() | Out-Null
$ws.Range("A114:A115").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("E113").Copy() | Out-Null
$ws.Range("E114:E115").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Row 114 ---
$ws.Range("A114").Value = 112
$ws.Range("B114").Value = 7749770
$ws.Range("C114").Value = "India Super League"
$ws.Range("D114").Value = "India Super League"
$ws.Range("E114").Value = 45381.35416666666
$ws.Range("F114").Value = "Bengaluru"
$ws.Range("G114").Value = "Odisha FC"
$ws.Range("K114").Value = 2.55
$ws.Range("L114").Value = 3.3
$ws.Range("M114").Value = 2.55
$ws.Range("N114").Value = 2.55
$ws.Range("O114").Value = 3.3
$ws.Range("P114").Value = 2.55
$ws.Range("Q114").Value = 0
$ws.Range("R114").Value = 1.9
$ws.Range("S114").Value = 1.9
$ws.Range("T114").Value = 2.5
$ws.Range("U114").Value = 1.8
$ws.Range("V114").Value = 2
$ws.Range("W114").Value = 0
$ws.Range("X114").Value = 0
$ws.Range("Y114").Value = 0
$ws.Range("Z114").Value = 0
$ws.Range("AA114").Value = 0

# --- Row 115 ---
$ws.Range("A115").Value = 113
$ws.Range("B115").Value = 7749469
$ws.Range("C115").Value = "India Super League"
$ws.Range("D115").Value = "India Super League"
$ws.Range("E115").Value = 45381.45833333334
$ws.Range("F115").Value = "Jamshedpur FC"
$ws.Range("G115").Value = "Kerala Blasters"
$ws.Range("K115").Value = 2.2
$ws.Range("L115").Value = 3.25
$ws.Range("M115").Value = 3
$ws.Range("N115").Value = 2.2
$ws.Range("O115").Value = 3.25
$ws.Range("P115").Value = 3
$ws.Range("Q115").Value = -0.25
$ws.Range("R115").Value = 1.975
$ws.Range("S115").Value = 1.825
$ws.Range("T115").Value = 2.5
$ws.Range("U115").Value = 1.925
$ws.Range("V115").Value = 1.875
$ws.Range("W115").Value = 0
$ws.Range("X115").Value = 0
$ws.Range("Y115").Value = 0
$ws.Range("Z115").Value = 0
$ws.Range("AA115").Value = 0
